$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '27.467.36'
$ws.Cells.Item(2, 5).Value2 = '  -0.80%  '
$ws.Cells.Item(3, 4).Value2 = '1.569.36'
$ws.Cells.Item(3, 5).Value2 = '  -0.89%  '
$ws.Cells.Item(4, 5).Value2 = '  -0.27%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = '208.86'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value2 = '  +1.00%  '
$ws.Cells.Item(6, 5).Value2 = '  -0.91%  '
$ws.Cells.Item(7, 5).Value2 = '  -0.25%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = '22.19'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value2 = '  -0.19%  '
$ws.Cells.Item(9, 5).Value2 = '  -1.13%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value2 = '0.0592'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value2 = '  +0.39%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = '0.0867'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value2 = '  -0.13%  '
$ws.Cells.Item(12, 4).Value2 = '1.792.93'
$ws.Cells.Item(12, 5).Value2 = '  -0.88%  '
$ws.Cells.Item(13, 4).Value2 = '1.560.33'
$ws.Cells.Item(13, 5).Value2 = '  -1.22%  '
$ws.Cells.Item(14, 5).Value2 = '  -0.96%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value2 = '0.519'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = '63.80'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value2 = '  +0.92%  '
$ws.Cells.Item(17, 4).Value2 = '27.471.45'
$ws.Cells.Item(17, 5).Value2 = '  -0.68%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = '214.48'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value2 = '  -2.02%  '
$ws.Cells.Item(19, 5).Value2 = '  +0.01%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = '7.28'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value2 = '  -0.55%  '
$ws.Cells.Item(21, 5).Value2 = '  -0.28%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = '4.12'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value2 = '  -0.29%  '
$ws.Cells.Item(23, 5).Value2 = '  +0.84%  '
$ws.Cells.Item(24, 5).Value2 = '  +1.82%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = '152.83'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value2 = '  -1.19%  '
$ws.Cells.Item(26, 5).Value2 = '  -0.27%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = '6.71'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value2 = '  -2.28%  '
$ws.Cells.Item(28, 5).Value2 = '  -0.52%  '
$ws.Cells.Item(29, 5).Value2 = '  -1.72%  '
$ws.Cells.Item(30, 5).Value2 = '  -0.05%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = '3.20'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value2 = '  -0.92%  '
$ws.Cells.Item(33, 4).Value2 = '1.376.77'
$ws.Cells.Item(33, 5).Value2 = '  -0.16%  '
$ws.Cells.Item(34, 5).Value2 = '  +1.90%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = '1.54'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value2 = '  +1.77%  '
$ws.Cells.Item(36, 5).Value2 = '  -0.84%  '
$ws.Cells.Item(37, 5).Value2 = '  -1.93%  '
$ws.Cells.Item(38, 5).Value2 = '  +1.41%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = '0.544'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value2 = '  +1.44%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = '0.828'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value2 = '  +1.07%  '
$ws.Cells.Item(41, 5).Value2 = '  -0.24%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = '0.980'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value2 = '  +0.48%  '
$ws.Cells.Item(43, 5).Value2 = '  +3.82%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = '64.30'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value2 = '  +1.40%  '
$ws.Cells.Item(45, 2).Value2 = 'MXToken'
$ws.Cells.Item(45, 3).Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = '2.17'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value2 = '  -0.15%  '
$ws.Cells.Item(46, 2).Value2 = 'FraxShare'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = '5.27'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value2 = '  +1.09%  '
$ws.Cells.Item(47, 4).Value2 = '1.704.87'
$ws.Cells.Item(47, 5).Value2 = '  -0.90%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = '85.32'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value2 = '  -3.37%  '
$ws.Cells.Item(49, 4).Value2 = '0.0₇0998'
$ws.Cells.Item(49, 5).Value2 = '  -0.75%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = '0.0959'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value2 = '  -1.35%  '
